# Fruta / hortaliza, semanal
# Insert a new daily price record for Vega Modelo de Temuco - Zanahoria as
# row 330, pushing the existing rows 330:405 down to 331:406.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row before the current row 330 (shifts 330..405 -> 331..406)
$ws.Rows.Item(330).Insert()

# Populate the newly inserted row 330 with the new record
$ws.Range("A330").Value = 10
$ws.Range("B330").Value = "Vega Modelo de Temuco"
$ws.Range("C330").Value = "La Araucanía"
$ws.Range("D330").Value = 44943
$ws.Range("E330").Value = 9
$ws.Range("F330").Value = 100114013
$ws.Range("G330").Value = "Zanahoria"
$ws.Range("H330").Value = "Sin especificar"
$ws.Range("I330").Value = "Primera"
$ws.Range("J330").Value = 115
$ws.Range("K330").Value = 10000
$ws.Range("L330").Value = 10000
$ws.Range("M330").Value = 10000
$ws.Range("N330").Value = "$/saco 20 kilos"
$ws.Range("O330").Value = "Región de La Araucanía"
$ws.Range("P330").Value = 500
$ws.Range("Q330").Value = 20
$ws.Range("R330").Value = "Hortaliza"

# Match the date-format style used by the rest of column D
$ws.Range("D330").NumberFormat = $ws.Range("D331").NumberFormat
